# HL update code for workflow 4
# Add a new "Manager" row for the existing Human Resources / nguyenltt
# contact (same email already used by row 5), with a matching mailto
# hyperlink, then reproduce the book's post-edit active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: department contact email (reuses the row-5 address) + role.
$ws.Range("A8").Value = "nguyenltt22411@st.uel.edu.vn"
$ws.Range("B8").Value = "Manager"

# Wire up the mailto: hyperlink on the new email cell, then restore the
# "Hyperlink" cell style (Hyperlinks.Add otherwise stamps its own style).
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:nguyenltt22411@st.uel.edu.vn")
$ws.Range("A8").Style = "Hyperlink"

# Match the saved selection state.
$ws.Range("B2").Select()
